# faturamento_diario.xlsx update
# - Corrige valores de total_venda para os dias 26 e 28 de maio/2025
# - Insere o dia 29 de maio/2025 (estava faltando), deslocando os dados
#   de abril/2025, março/2025 e fevereiro/2025 uma linha para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige valores existentes de maio/2025
$ws.Range("B27").Value = 27831.77
$ws.Range("B29").Value = 21718.63

# Insere uma nova linha na posicao 30 (antes dos dados de abril/2025),
# deslocando tudo a partir dali uma linha para baixo.
$ws.Rows.Item(30).Insert()

# Preenche a nova linha com o dia 29 de maio/2025
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 19445.75
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 2025
$ws.Range("E30").Value = "05/2025"

Write-Host "ok"
